$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A99").Value = "filter coffee mocha"
$ws.Range("B99").Value = "sweetened oat milk: 0.67 cup,sugar: 6.0 tsp,black coffee: 1.0 serving"
$ws.Range("C99").Value = 2
$ws.Range("D99").Value = 94.2
$ws.Range("E99").Value = 1.005
$ws.Range("F99").Value = 1.675
$ws.Range("G99").Value = 57.86

$ws.Range("A100").Value = "indian style vegan tofu pizza"
$ws.Range("B100").Value = "whole wheat pizza crust: 1.0,ragu pizza sauce: 0.5 cup,vegan mozzarella cheese: 1.0 cup,green capsicum: 0.5,mushroom: 0.5 cup,onion: 0.5,tomato: 1.5,firm tofu: 0.4 block,olive oil: 2.0 tsp"
$ws.Range("C100").Value = 2
$ws.Range("D100").Value = 424
$ws.Range("E100").Value = 25
$ws.Range("F100").Value = 14.5
$ws.Range("G100").Value = 54.75
